$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C3").Value = 174232
$ws.Range("C4").Value = 164226
$ws.Range("C7").Value = 5.74
$ws.Range("C8").Value = 64.48
